$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture the current row 5 / row 6 values (B:G never change format,
#         only their values get swapped between the two rows) ---
$b5 = $ws.Range("B5").Value()
$c5 = $ws.Range("C5").Value()
$d5 = $ws.Range("D5").Value()
$e5 = $ws.Range("E5").Value()
$f5 = $ws.Range("F5").Value()
$g5 = $ws.Range("G5").Value()

$a6 = $ws.Range("A6").Value()
$b6 = $ws.Range("B6").Value()
$c6 = $ws.Range("C6").Value()
$d6 = $ws.Range("D6").Value()
$e6 = $ws.Range("E6").Value()
$f6 = $ws.Range("F6").Value()
$g6 = $ws.Range("G6").Value()

# --- 2. Row 6's cylinder cell (A6) used the merged "6" label format
#         (vertical-top, border 5). Re-format it like the plain data rows
#         (A5 / A8) before the swap touches its value. ---
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 3. Swap the values: row 5 becomes the old row 6 data, row 6 becomes
#         the old row 5 data. ---
$ws.Range("A5").Value = $a6
$ws.Range("B5").Value = $b6
$ws.Range("C5").Value = $c6
$ws.Range("D5").Value = $d6
$ws.Range("E5").Value = $e6
$ws.Range("F5").Value = $f6
$ws.Range("G5").Value = $g6

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = $b5
$ws.Range("C6").Value = $c5
$ws.Range("D6").Value = $d5
$ws.Range("E6").Value = $e5
$ws.Range("F6").Value = $f5
$ws.Range("G6").Value = $g5

# --- 4. Row 6/7 no longer share a single merged "6" cylinder label -
#         unmerge, then give A7 its own value using the plain data-row
#         format (matching A5/A6/A8). ---
$ws.Range("A6:A7").UnMerge()
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A7").Value = 6
